$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GH")

# Row 4 (Inventory) updates
$ws.Range("B4").Value = 23000000.0
$ws.Range("C4").Value = 28000000.0
$ws.Range("D4").Value = 20000000.0
$ws.Range("E4").Value = 25000000.0
$ws.Range("F4").Value = 15000000.0

# Row 14 (Accounts Payable) updates
$ws.Range("B14").Value = 7000000.0
$ws.Range("C14").Value = 18000000.0
$ws.Range("D14").Value = 16000000.0
$ws.Range("E14").Value = 24000000.0
$ws.Range("F14").Value = 16000000.0
